# live_trading_results.xlsx - trade close + two new trade rows opened.
# Trade #92 closed at 2026-02-18 00:29:39 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's
# auto-detection turning date/time-looking strings (e.g. "2026-02-18") into
# serial-number dates, and without leaving a stray NumberFormat-derived
# style behind on the cell once we're done.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ===========================================================================
# Summary sheet
# ===========================================================================
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.09
$wsSummary.Range("B6").Value = 120
$wsSummary.Range("B9").Value = 48.33

# ===========================================================================
# Strategy Status sheet - MarketMaking row
# ===========================================================================
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D6").Value = 40
$wsStatus.Range("G6").Value = 50

# ===========================================================================
# All Trades sheet - close trade #120 (row 121) and append two new open
# trades (rows 150 and 151)
# ===========================================================================
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(121, 7).Value = 0.009792
Set-TextValue $wsAll.Cells.Item(121, 8) "CLOSED"
$wsAll.Cells.Item(121, 9).Value = -2.076
$wsAll.Cells.Item(121, 10).Value = -0
$wsAll.Cells.Item(121, 11).Value = 99.54000000000001
Set-TextValue $wsAll.Cells.Item(121, 12) "early_exit"
$wsAll.Cells.Item(121, 13).Value = 0.13

# New row 150: trade #149, momentum strategy, still open
$wsAll.Cells.Item(150, 1).Value = 149
Set-TextValue $wsAll.Cells.Item(150, 2) "2026-02-18"
Set-TextValue $wsAll.Cells.Item(150, 3) "00:29:32"
Set-TextValue $wsAll.Cells.Item(150, 4) "momentum"
Set-TextValue $wsAll.Cells.Item(150, 5) "UP"
$wsAll.Cells.Item(150, 6).Value = 0.01
Set-TextValue $wsAll.Cells.Item(150, 8) "OPEN"
$wsAll.Cells.Item(150, 9).Value = 0
$wsAll.Cells.Item(150, 10).Value = 0
$wsAll.Cells.Item(150, 11).Value = 99.23374292899115
$wsAll.Cells.Item(150, 13).Value = 0
$wsAll.Cells.Item(150, 14).Value = 0
$wsAll.Cells.Item(150, 15).Value = 0
$wsAll.Cells.Item(150, 16).Value = 0.9
Set-TextValue $wsAll.Cells.Item(150, 17) "Upward momentum: 2.941% over 10 samples"

# New row 151: trade #150, MarketMaking strategy, still open
$wsAll.Cells.Item(151, 1).Value = 150
Set-TextValue $wsAll.Cells.Item(151, 2) "2026-02-18"
Set-TextValue $wsAll.Cells.Item(151, 3) "00:29:33"
Set-TextValue $wsAll.Cells.Item(151, 4) "MarketMaking"
Set-TextValue $wsAll.Cells.Item(151, 5) "UP"
$wsAll.Cells.Item(151, 6).Value = 0.01
Set-TextValue $wsAll.Cells.Item(151, 8) "OPEN"
$wsAll.Cells.Item(151, 9).Value = 0
$wsAll.Cells.Item(151, 10).Value = 0
$wsAll.Cells.Item(151, 11).Value = 99.53967800952272
$wsAll.Cells.Item(151, 13).Value = 0
$wsAll.Cells.Item(151, 14).Value = 0
$wsAll.Cells.Item(151, 15).Value = 0
$wsAll.Cells.Item(151, 16).Value = 0.6
Set-TextValue $wsAll.Cells.Item(151, 17) "Normal spread capture: 190 bps"

# ===========================================================================
# momentum sheet - mirror of the new trade #149 row (row 36)
# ===========================================================================
$wsMomentum = $wb.Worksheets.Item("momentum")

$wsMomentum.Cells.Item(36, 1).Value = 149
Set-TextValue $wsMomentum.Cells.Item(36, 2) "2026-02-18"
Set-TextValue $wsMomentum.Cells.Item(36, 3) "00:29:32"
Set-TextValue $wsMomentum.Cells.Item(36, 4) "momentum"
Set-TextValue $wsMomentum.Cells.Item(36, 5) "UP"
$wsMomentum.Cells.Item(36, 6).Value = 0.01
Set-TextValue $wsMomentum.Cells.Item(36, 8) "OPEN"
$wsMomentum.Cells.Item(36, 9).Value = 0
$wsMomentum.Cells.Item(36, 10).Value = 0
$wsMomentum.Cells.Item(36, 11).Value = 99.23374292899115
$wsMomentum.Cells.Item(36, 12).Value = 0
$wsMomentum.Cells.Item(36, 13).Value = 0
$wsMomentum.Cells.Item(36, 14).Value = 0.9
Set-TextValue $wsMomentum.Cells.Item(36, 15) "Upward momentum: 2.941% over 10 samples"
$wsMomentum.Cells.Item(36, 17).Value = 0

# ===========================================================================
# MarketMaking sheet - close trade #120 (row 41) and append new trade #150
# (row 62)
# ===========================================================================
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Cells.Item(41, 7).Value = 0.009792
Set-TextValue $wsMM.Cells.Item(41, 8) "CLOSED"
$wsMM.Cells.Item(41, 9).Value = -2.076
$wsMM.Cells.Item(41, 10).Value = -0
$wsMM.Cells.Item(41, 11).Value = 99.54000000000001
Set-TextValue $wsMM.Cells.Item(41, 16) "early_exit"
$wsMM.Cells.Item(41, 17).Value = 0.13

# New row 62: trade #150, MarketMaking strategy, still open
$wsMM.Cells.Item(62, 1).Value = 150
Set-TextValue $wsMM.Cells.Item(62, 2) "2026-02-18"
Set-TextValue $wsMM.Cells.Item(62, 3) "00:29:33"
Set-TextValue $wsMM.Cells.Item(62, 4) "MarketMaking"
Set-TextValue $wsMM.Cells.Item(62, 5) "UP"
$wsMM.Cells.Item(62, 6).Value = 0.01
Set-TextValue $wsMM.Cells.Item(62, 8) "OPEN"
$wsMM.Cells.Item(62, 9).Value = 0
$wsMM.Cells.Item(62, 10).Value = 0
$wsMM.Cells.Item(62, 11).Value = 99.53967800952272
$wsMM.Cells.Item(62, 12).Value = 0
$wsMM.Cells.Item(62, 13).Value = 0
$wsMM.Cells.Item(62, 14).Value = 0.6
Set-TextValue $wsMM.Cells.Item(62, 15) "Normal spread capture: 190 bps"
$wsMM.Cells.Item(62, 17).Value = 0
